$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 69: Course, Hours, Notes for "Finish 3 hard problems from 3.6"
$ws.Range("B69").Value = "JS101: Programming Foundations with JavaScript"
$ws.Range("C69").Value = 0.75
$ws.Range("D69").Value = "3 questions from 3.6"

# Update selection to match the new active cell
$ws.Range("D69").Select()
